# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are plain text in the sheet (note the thousands-dot
# notation, e.g. "34.114.19"), but several of them also happen to look like
# ordinary decimals (e.g. "226.41"). Excel auto-converts a bare decimal-looking
# string assigned via .Value into a real number, which would corrupt values like
# "1.80" (trailing zero lost) and flip the stored cell type from text to number.
# Forcing NumberFormat to "@" (Text) before the write keeps it text, and resetting
# the style back to "Normal" afterwards avoids leaving a stray number-format on the
# cell (matching the original, unstyled inline-string cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Column D: Price
Set-TextValue $ws.Range("D2") "34.114.19"
Set-TextValue $ws.Range("D3") "1.787.72"
Set-TextValue $ws.Range("D5") "226.41"
Set-TextValue $ws.Range("D6") "0.546"
Set-TextValue $ws.Range("D8") "31.81"
Set-TextValue $ws.Range("D11") "0.0945"
Set-TextValue $ws.Range("D12") "2.046.23"
Set-TextValue $ws.Range("D13") "11.14"
Set-TextValue $ws.Range("D14") "1.776.05"
Set-TextValue $ws.Range("D15") "34.038.68"
Set-TextValue $ws.Range("D18") "67.98"
Set-TextValue $ws.Range("D19") "245.23"
Set-TextValue $ws.Range("D22") "10.84"
Set-TextValue $ws.Range("D25") "161.56"
Set-TextValue $ws.Range("D26") "7.14"
Set-TextValue $ws.Range("D27") "16.29"
Set-TextValue $ws.Range("D34") "1.80"
Set-TextValue $ws.Range("D35") "1.459.54"
Set-TextValue $ws.Range("D36") "2.42"
Set-TextValue $ws.Range("D37") "0.645"
Set-TextValue $ws.Range("D39") "1.03"
Set-TextValue $ws.Range("D40") "80.19"
Set-TextValue $ws.Range("D44") "13.45"
Set-TextValue $ws.Range("D48") "0.0₆0135"
Set-TextValue $ws.Range("D49") "1.947.41"
Set-TextValue $ws.Range("D50") "106.16"

# Column E: Volume(1h) -- percentage text cells (e.g. "  +0.46%  "); these never
# parse as numbers so a plain assignment already keeps them as text.
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  +4.86%  "
$ws.Range("E36").Value = "  +9.75%  "
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("E51").Value = "  +0.01%  "
